# data_dictionary.xlsx -- "wrangle" / "explore" update
#
# Summary of the change (from the commit message and the target diff):
#   * Several data-dictionary rows described fields whose data_type was
#     "bool"; those special encoding fields are now represented as
#     "int64" instead, so every occurrence of "bool" in column H
#     (data_type) becomes "int64". Because no cell references the
#     string "bool" any more, it naturally drops out of the shared
#     string table.
#   * The one-hot / helper encoding columns for company size, number of
#     researchers and ideal-conference-size (rows 81-90) are now
#     produced by the wrangle step, so column E (in_wrangle) flips from
#     FALSE to TRUE for those rows.
#   * The ideal-conference-size qids are renumbered so they sort/line
#     up correctly:
#       q17tws -> q17cs1   (row 87 - "workshop or retreat")
#       q17cs0 -> q17cs2   (row 88 - "fewer than 300 attendees")
#       q17cs1 -> q17cs3   (row 89 - "300-500 attendees")
#       q17cs4 stays       (row 90 - "more than 500 attendees")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuestionList")

# --- 1) data_type: bool -> int64 (column H) ---------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Value() -eq "bool") {
        $cell.Value = "int64"
    }
}

# --- 2) in_wrangle -> TRUE for the new encoding helper rows (column E) -----
for ($r = 81; $r -le 90; $r++) {
    $ws.Cells.Item($r, 5).Value = $true
}

# --- 3) Renumber the ideal-conference-size qids (column A) -----------------
$ws.Cells.Item(87, 1).Value = "q17cs1"
$ws.Cells.Item(88, 1).Value = "q17cs2"
$ws.Cells.Item(89, 1).Value = "q17cs3"

# --- 4) Scroll the frozen-header view down toward the bottom of the sheet --
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$ws.Range("A90").Select() | Out-Null
$win.ScrollRow = 75 | Out-Null
